$d = $word.ActiveDocument

# Replace the bookmark-style placeholder text in the first paragraph,
# also absorbing the trailing space run so only a single run remains.
$d.Content.Find.Execute("**ID__AFFARS_5316_topic_23__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_5316_603_2__ID**", 2)

# Update the first paragraph's formatting: add a paragraph border (space-only,
# no visible line) and widen the left indent.
$p1 = $d.Paragraphs(1)
$b = $p1.Borders
$b.DistanceFromTop = 5
$b.DistanceFromLeft = 5
$b.DistanceFromBottom = 5
$b.DistanceFromRight = 5
$p1.Range.ParagraphFormat.LeftIndent = 11.25
